$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 20
    3  = 20
    4  = 20
    5  = 70
    6  = 20
    7  = 20
    8  = 20
    9  = 70
    10 = 25
    11 = 25
    12 = 25
    13 = 70
    14 = 25
    15 = 25
    16 = 25
    17 = 70
    18 = 25
    19 = 25
    20 = 25
    21 = 70
    22 = 25
    23 = 25
    24 = 25
    25 = 70
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row]
}
